$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text: " brand" (leading space) -> "brand"
$ws.Range("D1").Value = "brand"

# Update the active selection/cell as recorded by the editor
$ws.Range("F8").Select()
